$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "${Date}" -> "${currentDate}" (with proofErr spell-check wrapping, as
#    Word does for camelCase merge-field names it does not recognise).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$dateFound = $r1.Find.Execute('${Date}', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($dateFound) {
    $range = $r1.Paragraphs(1).Range
    $xml = '<w:p ' + $wNs + ' w:rsidR="00495823" w:rsidRDefault="00B751D8"><w:pPr><w:spacing w:after="0" w:line="100" w:lineRule="atLeast"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>currentDate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>}</w:t></w:r></w:p>'
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 2) "Dear ${Employee Name}," -> "Dear ${personFirstName},"
# ---------------------------------------------------------------------------
$r2 = $d.Content
$dearFound = $r2.Find.Execute('Dear ${Employee Name},', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($dearFound) {
    $range = $r2.Paragraphs(1).Range
    $xml = '<w:p ' + $wNs + ' w:rsidR="00495823" w:rsidRDefault="00B751D8"><w:pPr><w:spacing w:after="0" w:line="100" w:lineRule="atLeast"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Dear </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>personFirstName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>},</w:t></w:r></w:p>'
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3) The "This letter is here to inform you, ${Employee Name}, ... case
#    file "${Case Number}: ${Case Title}"." paragraph ->
#    "This letter is here to inform you, ${personFirstName}, ... case file
#    "${caseNumber}: ${caseTitle}"."
# ---------------------------------------------------------------------------
$r3 = $d.Content
$letterFound = $r3.Find.Execute('This letter is here to inform you,', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($letterFound) {
    $range = $r3.Paragraphs(1).Range
    $xml = '<w:p ' + $wNs + ' w:rsidR="00495823" w:rsidRDefault="00B751D8"><w:pPr><w:spacing w:after="0" w:line="100" w:lineRule="atLeast"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">This letter is here to inform you, </w:t></w:r>' + `
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>${</w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>personFirstName</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>},</w:t></w:r>' + `
        '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> is under investigation currently and you have been requested to be interviewed as part of case file </w:t></w:r>' + `
        '<w:r w:rsidR="00F40982" w:rsidRPr="00F40982"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&quot;</w:t></w:r>' + `
        '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>${</w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>c</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>aseNumber</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>}: ${</w:t></w:r>' + `
        '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>c</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
        '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>aseTitle}</w:t></w:r>' + `
        '<w:r w:rsidR="00F40982" w:rsidRPr="00F40982"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&quot;</w:t></w:r>' + `
        '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r>' + `
        '</w:p>'
    $range.InsertXML($xml)
}
